# Commit: feat(main): add lab-4 files
#
# Appends a new "Выводы" (Conclusions) section at the end of the document:
#   - a Heading2 paragraph with the title "Выводы"
#   - a FirstParagraph-styled paragraph containing four bold "claim" phrases
#     each followed by a plain-text explanatory sentence
# and wraps that whole new section in a bookmark named "выводы" (mirroring
# the bookmarkStart/bookmarkEnd pair added around it in the target revision).

$d = $word.ActiveDocument

# Remember where the new content begins: right after the current last
# paragraph ("Выполнение lab5-4"), i.e. at the end of the document body.
$rng = $d.Content
$rng.Collapse(0)
$insertStart = $rng.Start

# Build the new "Выводы" heading + body paragraph as a literal OOXML
# fragment so that run-level formatting (bold -> <w:b/><w:bCs/>) and
# paragraph styles come out exactly as in the target revision.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t xml:space="preserve">Выводы</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Оба подхода успешно решают поставленную задачу</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- программы lab5-3.asm и lab5-4.asm корректно выполняют ввод и вывод строки</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Использование внешних подпрограмм упрощает разработку</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- подход с in_out.asm требует меньше кода и более читаем</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Прямые системные вызовы дают больше контроля над процессом</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- ручное управление системными вызовами обеспечивает лучший контроль</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Программы корректно работают с русскими символами</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- обеспечена правильная обработка кириллических символов</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $rng.InsertXML($xml)

# Wrap the freshly-inserted heading + paragraph in a new bookmark called
# "выводы", matching the bookmarkStart/bookmarkEnd pair introduced around
# this section in the diff.
$insertEnd = $d.Content.End
$newSectionRange = $d.Range($insertStart, $insertEnd)
$d.Bookmarks.Add("выводы", $newSectionRange)
